$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New record row ("MCH206-1" / Series / 1 Box / location note)
$ws.Range("A2").Value = "MCH206-1"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 24D | GRAP COUNT NUMER: NONE"

# Match the row's font (Calibri 10pt, text1 theme color) like the rest of the sheet body
$ws.Range("A2").Font.Name = "Calibri"
$ws.Range("A2").Font.Size = 10
$ws.Range("A2").Font.ThemeColor = 1

# Reuse that same resolved style for the rest of the row instead of re-deriving it
$ws.Range("A2").Copy() | Out-Null
$ws.Range("C2:H2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Keep the header row frozen and select the new row, same as the saved view
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A2:J2").Select() | Out-Null
